$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting rows 11:26 down to 12:27.
# This mirrors Excel's native "Insert Sheet Rows" action (copies formatting
# from the row above by default).
$ws.Rows(11).Insert()

# Populate the new row 11 with this week's data point.
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(11, 3).Value = "Los Lagos"
$ws.Cells.Item(11, 4).Value = 44838
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 100112035
$ws.Cells.Item(11, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 22000
$ws.Cells.Item(11, 12).Value = 22000
$ws.Cells.Item(11, 13).Value = 22000
$ws.Cells.Item(11, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(11, 16).Value = 1467
$ws.Cells.Item(11, 17).Value = 15
$ws.Cells.Item(11, 18).Value = "Hortaliza"
